# Implement the "Upload Lessons Learned, View Lessons Learned and View Action
# Items of Add Project in Dashboard" change: add a new "addProject" worksheet
# at the end of the workbook (becoming the active sheet) with a 2x4 table of
# shared-string-backed titles/labels.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (and becomes the active sheet, like in the diff).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "addProject"

# Header row (row 1) / value row (row 2)
$newSheet.Range("A1").Value = "uploadLessonsLearnedTitle"
$newSheet.Range("B1").Value = "viewLessonsLearnedTitle"
$newSheet.Range("C1").Value = "uploadActionItemsTitle"
$newSheet.Range("D1").Value = "viewActionItemsTitle"

$newSheet.Range("A2").Value = "Upload Lessons Learned"
$newSheet.Range("B2").Value = "View Lessons Learned"
$newSheet.Range("C2").Value = "Upload Action Items"
$newSheet.Range("D2").Value = "View Action items"

# Column widths matching the authored workbook.
$newSheet.Columns.Item(1).ColumnWidth = 25.833333333333332
$newSheet.Columns.Item(2).ColumnWidth = 26.833333333333332
$newSheet.Columns.Item(3).ColumnWidth = 26.0
$newSheet.Columns.Item(4).ColumnWidth = 26.0

# Leave the selection on D2, matching the authored workbook's saved state.
[void]$newSheet.Range("D2").Select()
